# Updated test data for normal load, cable capacitance etc
#
# Workbook has 4 sheets:
#   1 "Add Devices Loop A"
#   2 "Add Devices Loop B"
#   3 "Add Devices Loop C"
#   4 "Add Devices Loop D"

$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item(1)
$wsB = $wb.Worksheets.Item(2)
$wsC = $wb.Worksheets.Item(3)
$wsD = $wb.Worksheets.Item(4)

# ---- Loop A : updated normal-load / volt-drop figures ----
$wsA.Range("E7").Value = 279
$wsA.Range("F7").Value = 0.24
$wsA.Range("G7").Value = 0.42

# ---- Loop B : updated normal-load / volt-drop figures ----
$wsB.Range("E6").Value = 337
$wsB.Range("F6").Value = 0.24
$wsB.Range("G6").Value = 0.41

# ---- Loop C : no data changes ----

# ---- Loop D : updated cable-capacitance / volt-drop figures ----
# F6 used to hold the quoted-text "0.10" - replace with the real number
$wsD.Range("F6").Value = 0.03
$wsD.Range("G6").Value = 0.06
$wsD.Range("F7").Value = 0.07
$wsD.Range("G7").Value = 0.14

# ---- Selections: touch the sheets that only move their cursor first ----
$wsB.Range("B5").Select() | Out-Null
$wsD.Range("F6").Select() | Out-Null

# Loop A becomes the active tab/selection last, matching the saved view
$wsA.Range("G8").Select() | Out-Null
